$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 6 block (rows 27-29)
$ws.Range("C27").Value = 6936
$ws.Range("C28").Value = 1700
$ws.Range("C29").Value = 1700

# Day 7 block (rows 33-35)
$ws.Range("C33").Value = 6936
$ws.Range("C34").Value = 1775
$ws.Range("C35").Value = 1775

# Day 8 block (rows 39-41)
$ws.Range("C39").Value = 6936
$ws.Range("C40").Value = 1906
$ws.Range("C41").Value = 1906

# Update sheet view: scroll position and active selection cell
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G39").Select()
